$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 9, pushing the existing rows 9-76 down to 10-77.
$ws.Rows(9).Insert()

# Populate the newly inserted row 9 with the new daily price entry.
$ws.Range("A9").Value = 10
$ws.Range("B9").Value = "Vega Modelo de Temuco"
$ws.Range("C9").Value = "La Araucanía"
$ws.Range("D9").Value = 45022
$ws.Range("E9").Value = 9
$ws.Range("F9").Value = 100112010
$ws.Range("G9").Value = "Achicoria"
$ws.Range("H9").Value = "Sin especificar"
$ws.Range("I9").Value = "Primera"
$ws.Range("J9").Value = 65
$ws.Range("K9").Value = 10000
$ws.Range("L9").Value = 10000
$ws.Range("M9").Value = 10000
$ws.Range("N9").Value = "$/caja 18 unidades"
$ws.Range("O9").Value = "Región Metropolitana"
$ws.Range("P9").Value = 556
$ws.Range("Q9").Value = 18
$ws.Range("R9").Value = "Hortaliza"
